$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Pearson Edexcel logo images (PNG) in the footers: image2.png -> image1.png
$f1 = $sec.Footers(1)   # wdHeaderFooterPrimary
$f2 = $sec.Footers(2)   # wdHeaderFooterFirstPage
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
    $f1.Range.InlineShapes(1).Name = "image1.png"
}
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
    $f2.Range.InlineShapes(1).Name = "image1.png"
}

# BTEC logo images (JPG) in the headers: image1.jpg -> image2.jpg
$h1 = $sec.Headers(1)   # wdHeaderFooterPrimary
$h2 = $sec.Headers(2)   # wdHeaderFooterFirstPage
if ($h1.Exists -and $h1.Range.InlineShapes.Count -ge 1) {
    $h1.Range.InlineShapes(1).Name = "image2.jpg"
}
if ($h2.Exists -and $h2.Range.InlineShapes.Count -ge 1) {
    $h2.Range.InlineShapes(1).Name = "image2.jpg"
}
